$d = $word.ActiveDocument

$newText = "Tämän oppaan kartat piirsi Jenik Hollan CzechGlobesta (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    if ($pr.Text -like "*Tämän oppaan kartat piirsi*" -and $pr.Text -like "*GaNight/2018*") {
        $start = $pr.Start
        $end = $pr.End - 1   # exclude the paragraph mark
        $target = $d.Range($start, $end)
        $target.Text = ""
        $insertion = $d.Range($start, $start)
        $insertion.InsertAfter($newText)
        break
    }
}
